$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap values of G2 and G3 (Actionneur Small column)
$g2 = $ws.Range("G2").Value2
$g3 = $ws.Range("G3").Value2
$ws.Range("G2").Value2 = $g3
$ws.Range("G3").Value2 = $g2

# Update the view: scroll so column B is the top-left visible column,
# and select cell G3
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Range("G3").Select()
